$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Mã nhóm" (group code) values to their shortened forms
$ws.Range("A2").Value = "TH1"
$ws.Range("A3").Value = "TH2"
$ws.Range("A4").Value = "LT"

# Update the active selection to match the saved cursor position
$ws.Range("D9").Select()
